# Update generated event-count figures ("想去人数") for the two sheets
# ("展览" and "全部类型") that list the 南宁·第五届小萌萌动漫嘉年华 row (F2)
# and the 南宁·草莓动漫节 row (F4), matching the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 313
    $ws.Range("F4").Value = 1262
}
